$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 6.5
$ws.Range("J7").Value = 2.82
$ws.Range("K7").Value = 1.98
$ws.Range("L7").Value = 3.85
$ws.Range("N7").Value = 6.7
$ws.Range("P7").Value = 2.55
$ws.Range("T7").Value = 2.42
$ws.Range("U7").Value = 1.83
$ws.Range("V7").Value = 1.78
$ws.Range("W7").Value = 6.5
$ws.Range("X7").Value = 10.25
$ws.Range("Y7").Value = 9
$ws.Range("AA7").Value = 20
$ws.Range("AG7").Value = 8.5
$ws.Range("AH7").Value = 17
$ws.Range("AI7").Value = 11.75
$ws.Range("AL7").Value = 45
$ws.Range("AO7").Value = 11.75
$ws.Range("AP7").Value = 20
$ws.Range("AR7").Value = 80
$ws.Range("AS7").Value = 250
$ws.Range("AT7").Value = 2.4
$ws.Range("AV7").Value = 60
$ws.Range("AY7").Value = 25
$ws.Range("AZ7").Value = 100
$ws.Range("G8").Value = 2.5
$ws.Range("I8").Value = 2.75
$ws.Range("J8").Value = 3.25
$ws.Range("L8").Value = 3.5
$ws.Range("X8").Value = 12
$ws.Range("Z8").Value = 26
$ws.Range("AA8").Value = 23
$ws.Range("AG8").Value = 8
$ws.Range("AH8").Value = 13
$ws.Range("AK8").Value = 23
$ws.Range("AL8").Value = 34
$ws.Range("AN8").Value = 4.5
$ws.Range("AR8").Value = 81
$ws.Range("G12").Value = 1.3
$ws.Range("H12").Value = 4.75
$ws.Range("I12").Value = 11
$ws.Range("K12").Value = 2.3
$ws.Range("Q12").Value = 2.03
$ws.Range("R12").Value = 1.83
$ws.Range("S12").Value = 1.4
$ws.Range("T12").Value = 2.75
$ws.Range("Y12").Value = 10
$ws.Range("Z12").Value = 7.5
$ws.Range("AC12").Value = 9
$ws.Range("AD12").Value = 9.5
$ws.Range("AG12").Value = 19
$ws.Range("AK12").Value = 81
$ws.Range("AP12").Value = 23
$ws.Range("AR12").Value = 51
$ws.Range("AT12").Value = 2.75
$ws.Range("AW12").Value = 10
$ws.Range("AZ12").Value = 301
$ws.Range("M13").Value = 1.06
$ws.Range("N13").Value = 8.970000000000001
$ws.Range("H14").Value = 2.95
$ws.Range("I14").Value = 2.42
$ws.Range("J14").Value = 3.4
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 3.05
$ws.Range("N14").Value = 7.2
$ws.Range("P14").Value = 2.67
$ws.Range("S14").Value = 1.42
$ws.Range("T14").Value = 2.47
$ws.Range("W14").Value = 8.5
$ws.Range("AA14").Value = 26
$ws.Range("AB14").Value = 35
$ws.Range("AC14").Value = 7.9
$ws.Range("AG14").Value = 7
$ws.Range("AI14").Value = 9.5
$ws.Range("AL14").Value = 35
$ws.Range("AO14").Value = 15.5
$ws.Range("AQ14").Value = 70
$ws.Range("AR14").Value = 100
$ws.Range("AS14").Value = 250
$ws.Range("AW14").Value = 4.3
$ws.Range("AX14").Value = 13
$ws.Range("AY14").Value = 21
$ws.Range("AZ14").Value = 60
$ws.Range("BA14").Value = 100
$ws.Range("BB14").Value = 300
